# Add a new "2022-Q4" sheet (fund-holdings detail) right before the
# existing "2022-Q3" sheet, and insert a corresponding summary row at the
# top of the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Totals sheet ("总计"): insert a new row 2 for 2022-Q4 and push the
#    previous rows down by one.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

$totals.Rows.Item(2).Insert()
# Re-apply the formatting of the row below (which carries the original
# template styling) onto the freshly inserted, blank row.
$totals.Range("A3:D3").Copy($totals.Range("A2:D2"))

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q4"
$totals.Cells.Item(2,3).Value = 9
$totals.Cells.Item(2,4).Value = 1.08

# Renumber the index column (A) for the rows that shifted down so the
# sequence stays 0,1,2,3,4,5.
$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(5,1).Value = 3
$totals.Cells.Item(6,1).Value = 4
$totals.Cells.Item(7,1).Value = 5

# ------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet itself, positioned right before the
#    "2022-Q3" sheet. Cloning the existing "2022-Q3" sheet gives us an
#    identical template (headers, styles, column layout) to start from.
# ------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$beforeSheet.Copy($beforeSheet)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template (old "2022-Q3") has 10 data rows (rows 2-11); the new
# quarter only needs 9 data rows (rows 2-10), so drop the extra one.
$q4.Rows.Item(11).Delete()

# Columns that hold numeric-looking text (fund code / size / position /
# weight / market value) must be forced to Text so Excel doesn't coerce
# them back into numbers.
$q4.Range("B2:B10").NumberFormat = "@"
$q4.Range("D2:G10").NumberFormat = "@"

function Set-FundRow($row, $idx, $code, $name, $size, $pos, $weight, $mv, $rank) {
    $q4.Cells.Item($row, 1).Value = $idx
    $q4.Cells.Item($row, 2).Value = $code
    $q4.Cells.Item($row, 3).Value = $name
    $q4.Cells.Item($row, 4).Value = $size
    $q4.Cells.Item($row, 5).Value = $pos
    $q4.Cells.Item($row, 6).Value = $weight
    $q4.Cells.Item($row, 7).Value = $mv
    $q4.Cells.Item($row, 8).Value = $rank
}

Set-FundRow 2 0 "011136" "广发盛兴混合A" "17.77" "87.23" "2.52" "0.4478" 10
Set-FundRow 3 1 "506007" "广发科创板两年定开混合" "5.31" "88.81" "5.69" "0.3021" 2
Set-FundRow 4 2 "012342" "广发瑞泽精选混合A" "5.19" "89.55" "2.79" "0.1448" 9
Set-FundRow 5 3 "013000" "广发盛泽一年持有期混合A" "2.47" "85.08" "2.67" "0.0659" 10
Set-FundRow 6 4 "011137" "广发盛兴混合C" "1.85" "87.23" "2.52" "0.0466" 10
Set-FundRow 7 5 "002133" "广发鑫益灵活配置混合" "1.33" "87.79" "3.06" "0.0407" 10
Set-FundRow 8 6 "005027" "光大保德信多策略优选一年定期开放灵活配置混合" "0.49" "53.72" "3.11" "0.0152" 8
Set-FundRow 9 7 "012343" "广发瑞泽精选混合C" "0.37" "89.55" "2.79" "0.0103" 9
Set-FundRow 10 8 "013001" "广发盛泽一年持有期混合C" "0.28" "85.08" "2.67" "0.0075" 10

Write-Output "2022-Q4 sheet added; worksheets now: $($wb.Worksheets.Count)"
